# Update res_line/pl_mw data for the 380 kV case (rows A=0..23, i.e. sheet rows 2-25).
# Only columns B, C, D, F, G, I, K, N change; E, H, J, L, M, O stay 0 and column A stays the index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 0.4487503271861613; "C2" = 0.08999107383239391; "D2" = 0.04770524831767631; "F2" = 1.11004860662316; "G2" = 0.002471323633778969; "I2" = 0.7769113166066148; "K2" = 0.517949952431934; "N2" = 1.689458135213776
    "B3" = 0.4097608860521973; "C3" = 0.08145825179261124; "D3" = 0.04674374881948395; "F3" = 1.107724365235526; "G3" = 0.002474448741029027; "I3" = 0.7815845219114372; "K3" = 0.4719564677101005; "N3" = 1.709139021880439
    "B4" = 0.3859817992284604; "C4" = 0.07625646039123524; "D4" = 0.04614765729770909; "F4" = 1.106997281224707; "G4" = 0.002476469716336441; "I4" = 0.7848870513603075; "K4" = 0.4439091284640142; "N4" = 1.721845837272049
    "B5" = 0.3763322768375872; "C5" = 0.07414605427430843; "D5" = 0.04590332244168138; "F5" = 1.106876814289443; "G5" = 0.002477319043971884; "I5" = 0.7863416716365172; "K5" = 0.4325282993694088; "N5" = 1.727180438054971
    "B6" = 0.3747324440066109; "C6" = 0.07379618756166906; "D6" = 0.04586266540986372; "F6" = 1.106867424545513; "G6" = 0.002477461632383036; "I6" = 0.7865897790572767; "K6" = 0.4306414679606121; "N6" = 1.728075691834567
    "B7" = 0.3858514974547234; "C7" = 0.07622796084149286; "D7" = 0.04614436784993714; "F7" = 1.106994944917759; "G7" = 0.002476481066480676; "I7" = 0.784906228446637; "K7" = 0.4437554452190113; "N7" = 1.721917148164312
    "B8" = 0.4352735853628644; "C8" = 0.08704117001462919; "D8" = 0.0473749220176245; "F8" = 1.109101798158598; "G8" = 0.002472380017462133; "I8" = 0.7784326681229388; "K8" = 0.5020514762730386; "N8" = 1.69611480788709
    "B9" = 0.5334578386647024; "C9" = 0.1085453898292599; "D9" = 0.04974194115971642; "F9" = 1.118798407244824; "G9" = 0.002465144794040314; "I9" = 0.7691803678266567; "K9" = 0.6178975618521463; "N9" = 1.650461044991346
    "B10" = 0.6063656493258236; "C10" = 0.1245327205656679; "D10" = 0.05145217378022693; "F10" = 1.129332850905811; "G10" = 0.002460315991601186; "I10" = 0.7644893292067181; "K10" = 0.7039470296913066; "N10" = 1.619937806600921
    "B11" = 0.6397011856194013; "C11" = 0.1318479307237226; "D11" = 0.0522238059102591; "F11" = 1.13486977353935; "G11" = 0.002458223908193135; "I11" = 0.7628144921491469; "K11" = 0.7432985594692809; "N11" = 1.606708376500631
    "B12" = 0.6523486933859317; "C12" = 0.1346241891272086; "D12" = 0.05251507316720705; "F12" = 1.137073830854462; "G12" = 0.00245744664621349; "I12" = 0.7622464394875337; "K12" = 0.7582296988695703; "N12" = 1.601793154531268
    "B13" = 0.6496237607783257; "C13" = 0.1340259982284238; "D13" = 0.05245238534139673; "F13" = 1.136594369486559; "G13" = 0.002457613378973705; "I13" = 0.7623658344186524; "K13" = 0.7550126993321271; "N13" = 1.602847530758684
    "B14" = 0.6407412231226886; "C14" = 0.1320762117660195; "D14" = 0.05224778744808134; "F14" = 1.135048949872655; "G14" = 0.002458159662786646; "I14" = 0.76276643108325; "K14" = 0.7445263619725893; "N14" = 1.606302103741294
    "B15" = 0.635303537327303; "C15" = 0.1308827128804353; "D15" = 0.05212234327362353; "F15" = 1.134116322654336; "G15" = 0.002458496224272836; "I15" = 0.7630204300323626; "K15" = 0.7381070216501939; "N15" = 1.6084304384453
    "B16" = 0.6041904771422253; "C16" = 0.1240555125755236; "D16" = 0.05140161622889394; "F16" = 1.12898600769411; "G16" = 0.002460454812012935; "I16" = 0.7646080354526461; "K16" = 0.701379464694952; "N16" = 1.620815591601596
    "B17" = 0.5851468002809384; "C17" = 0.1198781558372275; "D17" = 0.05095783114983732; "F17" = 1.126029652401542; "G17" = 0.002461683072361816; "I17" = 0.7656996817953043; "K17" = 0.678901205236258; "N17" = 1.628581561312814
    "B18" = 0.5742093441789109; "C18" = 0.1174794481927108; "D18" = 0.05070198004224125; "F18" = 1.124399319527882; "G18" = 0.002462399381381971; "I18" = 0.7663707754186788; "K18" = 0.6659918205464805; "N18" = 1.633110052375214
    "B19" = 0.5705088602204285; "C19" = 0.1166679724549908; "D19" = 0.05061525116487786; "F19" = 1.123859346290402; "G19" = 0.002462643604812289; "I19" = 0.7666054126271504; "K19" = 0.6616242834081163; "N19" = 1.634653915061572
    "B20" = 0.587172382221695; "C20" = 0.1203224286076932; "D20" = 0.051005134818773; "F20" = 1.126337106171221; "G20" = 0.002461551303546807; "I20" = 0.7655790012900141; "K20" = 0.681292036170845; "N20" = 1.627748472450108
    "B21" = 0.6433495890300378; "C21" = 0.1326487441912718; "D21" = 0.05230790827774712; "F21" = 1.135499961765603; "G21" = 0.00245799880059605; "I21" = 0.7626469690656421; "K21" = 0.7476056527423225; "N21" = 1.605284846402967
    "B22" = 0.6802047652643921; "C22" = 0.1407405509914383; "D22" = 0.05315389892761857; "F22" = 1.14211422242488; "G22" = 0.002455764226525133; "I22" = 0.7611164855318009; "K22" = 0.791117682029352; "N22" = 1.5911543335977
    "B23" = 0.6605217406475958; "C23" = 0.1364185092644448; "D23" = 0.05270288222202169; "F23" = 1.13852672290254; "G23" = 0.002456948906098228; "I23" = 0.7618979864084849; "K23" = 0.7678788067029245; "N23" = 1.598645598483998
    "B24" = 0.5862565823181285; "C24" = 0.1201215638266717; "D24" = 0.05098375105046671; "F24" = 1.126197890424905; "G24" = 0.002461610844516566; "I24" = 0.7656334254802388; "K24" = 0.6802110989224275; "N24" = 1.628124913409023
    "B25" = 0.5067607865452715; "C25" = 0.1026952851682381; "D25" = 0.04910660916768705; "F25" = 1.11557760885637; "G25" = 0.002467016246473719; "I25" = 0.7713140020346145; "K25" = 0.5863939546807444; "N25" = 1.662281823861226
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
